$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.659.78"
$ws.Range("E2").Value = "  +1.84%  "

$ws.Range("D3").Value = "3.162.69"
$ws.Range("E3").Value = "  +1.45%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'528.99"
$ws.Range("E5").Value = "  -0.38%  "

$ws.Range("D6").Value = "'140.15"
$ws.Range("E6").Value = "  +1.43%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").Value = "'0.535"
$ws.Range("E8").Value = "  +13.78%  "

$ws.Range("D9").Value = "'7.29"
$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("D10").Value = "'0.438"
$ws.Range("E10").Value = "  +6.61%  "

$ws.Range("D11").Value = "'0.112"
$ws.Range("E11").Value = "  +4.27%  "

$ws.Range("D12").Value = "'0.141"
$ws.Range("E12").Value = "  +2.57%  "

$ws.Range("D13").Value = "3.710.09"
$ws.Range("E13").Value = "  +1.66%  "

$ws.Range("D14").Value = "'25.75"
$ws.Range("E14").Value = "  +0.50%  "

$ws.Range("E15").Value = "  +3.50%  "

$ws.Range("D16").Value = "58.697.29"
$ws.Range("E16").Value = "  +1.72%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "'6.23"
$ws.Range("E17").Value = "  +3.20%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.149.06"
$ws.Range("E18").Value = "  +0.93%  "

$ws.Range("D19").Value = "'12.96"
$ws.Range("E19").Value = "  +2.23%  "

$ws.Range("D20").Value = "'8.12"
$ws.Range("E20").Value = "  +0.68%  "

$ws.Range("D21").Value = "'375.48"
$ws.Range("E21").Value = "  +4.07%  "

$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").Value = "'0.528"
$ws.Range("E23").Value = "  +4.83%  "

$ws.Range("D24").Value = "'69.69"
$ws.Range("E24").Value = "  +0.93%  "

$ws.Range("E25").Value = "  +0.62%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("D27").Value = "'8.33"
$ws.Range("E27").Value = "  +14.31%  "

$ws.Range("D28").Value = "0.0₃0860"
$ws.Range("E28").Value = "  -0.44%  "

$ws.Range("D29").Value = "'22.46"
$ws.Range("E29").Value = "  +5.42%  "

$ws.Range("E30").Value = "  +0.89%  "

$ws.Range("E31").Value = "  -1.28%  "

$ws.Range("D32").Value = "'5.14"
$ws.Range("E32").Value = "  +0.14%  "

$ws.Range("E33").Value = "  +0.87%  "

$ws.Range("D34").Value = "'6.31"
$ws.Range("E34").Value = "  +4.18%  "

$ws.Range("D35").Value = "'156.83"
$ws.Range("E35").Value = "  -1.57%  "

$ws.Range("E36").Value = "  +3.74%  "

$ws.Range("D37").Value = "2.697.01"
$ws.Range("E37").Value = "  +6.57%  "

$ws.Range("D38").Value = "'25.02"
$ws.Range("E38").Value = "  -1.73%  "

$ws.Range("E39").Value = "  +1.01%  "

$ws.Range("E40").Value = "  +3.25%  "

$ws.Range("E41").Value = "  +6.83%  "

$ws.Range("D42").Value = "'0.722"
$ws.Range("E42").Value = "  +3.56%  "

$ws.Range("D43").Value = "'39.17"
$ws.Range("E43").Value = "  +3.61%  "

$ws.Range("E44").Value = "  +7.52%  "

$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("D46").Value = "3.204.91"
$ws.Range("E46").Value = "  +1.45%  "

$ws.Range("E47").Value = "  +12.71%  "

$ws.Range("E48").Value = "  +1.83%  "

$ws.Range("D49").Value = "'0.979"
$ws.Range("E49").Value = "  +0.49%  "

$ws.Range("D50").Value = "'20.06"
$ws.Range("E50").Value = "  +1.75%  "

$ws.Range("D51").Value = "'0.748"
$ws.Range("E51").Value = "  +0.89%  "

